# Updated symbol list with GitHub Actions
# Refresh the Price (column D) and Volume(1h) (column E) figures for the
# cryptocurrency rows. Values are kept as plain text (matching the source
# data, which stores prices/percentages as strings, e.g. "0.01100").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "304.16"
    "D3" = "35.81"
    "E3" = "-4.31%"
    "D4" = "5.094"
    "E4" = "1.99%"
    "D5" = "0.07849"
    "E5" = "0.29%"
    "D6" = "2.113"
    "E6" = "-3.62%"
    "D7" = "7.923"
    "E7" = "-1.40%"
    "D8" = "0.9192"
    "E8" = "1.14%"
    "E9" = "0.90%"
    "D10" = "0.1858"
    "E10" = "-1.98%"
    "D11" = "0.08604"
    "E11" = "1.38%"
    "D12" = "0.03548"
    "E12" = "0.59%"
    "D13" = "0.09952"
    "E13" = "-0.13%"
    "D14" = "0.001441"
    "E14" = "-2.86%"
    "D15" = "0.005700"
    "E15" = "0.95%"
    "D16" = "3.468"
    "E16" = "-0.03%"
    "D17" = "4.101"
    "E17" = "1.66%"
    "E18" = "21.31%"
    "D19" = "0.3421"
    "E19" = "-1.21%"
    "D20" = "5.232"
    "E20" = "9.76%"
    "D21" = "0.1310"
    "E21" = "1.37%"
    "D23" = "0.04553"
    "E23" = "-1.78%"
    "D24" = "0.005065"
    "E24" = "5.38%"
    "D25" = "0.001234"
    "E25" = "0.37%"
    "D27" = "0.0004754"
    "D39" = "0.01838"
    "E39" = "4.61%"
    "D40" = "0.04735"
    "E40" = "0.28%"
    "D41" = "0.007575"
    "E41" = "-6.02%"
    "D42" = "0.1399"
    "E42" = "0.24%"
    "D43" = "0.007729"
    "E43" = "0.89%"
    "D44" = "0.002207"
    "E44" = "1.62%"
    "D45" = "0.01128"
    "E45" = "8.47%"
    "D46" = "0.00006328"
    "E46" = "4.48%"
    "D47" = "0.00000000751"
    "E47" = "0.08%"
    "D48" = "0.0005806"
    "E48" = "0.09%"
    "D49" = "47.10"
    "E49" = "676.48%"
    "D50" = "0.002002"
    "E50" = "-25.59%"
    "D51" = "0.00002102"
    "E51" = "0.08%"
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    # Force text storage so values like "0.01100" or "-4.31%" keep their
    # exact literal formatting instead of being parsed into numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
}
